$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[3] = @("Otto Sverdrup",350,60,40,0,50,10,1000,40,20,20,0,0,0,0,0,0,0,1200,5,5,90,0,0,1600,5,5,90,0,0,1200,10,300,30,0,0,10,0,350,60,40,0,50,10,300,10,0,90,10,30,20,40,0,0,0,100,30,10,60,0,0,100,10,30,60,0,0,50,30,50,30,0,0,10,0)
$data[4] = @("Tom Crean",350,60,40,0,100,50,1000,20,80,0,0,0,0,0,0,0,0,1200,5,5,90,0,0,600,5,5,90,0,0,600,10,600,30,0,0,10,0,350,60,40,0,50,10,300,20,80,0,0,0,0,0,0,0,0,100,30,10,60,0,0,500,10,30,60,0,0,500,30,50,30,0,0,10,0)
$data[5] = @("Helen Thayer",250,40,60,0,50,10,600,20,80,0,0,0,0,0,0,0,0,600,5,5,90,0,0,1600,5,5,90,0,0,1200,10,100,30,0,0,10,0,350,60,40,0,50,10,300,40,60,0,60,30,20,40,0,0,0,500,30,10,60,0,0,100,10,30,60,0,0,50,30,500,30,0,0,10,0)

foreach ($r in 3..5) {
    $rowVals = $data[$r]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

$ws.Range("A3:A5").Font.Bold = $true

$null = $ws.Range("A5").Select()
